$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "A Lag" row (row 2): column B = A coefficient, column C = C coefficient
$ws.Range("B2").Value = "-0.372***"
$ws.Range("C2").Value = "0.01*"

# Update the "C Lag" row (row 3): column B = A coefficient, column C = C coefficient
$ws.Range("B3").Value = "-3.464***"
$ws.Range("C3").Value = "-0.808***"
